# Applies the cryptos-list price/volume refresh described by the commit.
#
# Target cells fall into two buckets that need different handling:
#
#   1) Text that Excel will not try to reinterpret as a number (coin names,
#      URLs, multi-dot "thousand-separated" prices like "65.424.25", and the
#      "  +5.34%  " volume strings with their padding spaces/percent sign).
#      A plain Range.Value assignment is safe for these and leaves the cell
#      format untouched.
#
#   2) Plain decimal-looking prices ("418.27", "0.999", "1.00", ...). Excel
#      COM auto-coerces a bare numeric string to a real number on assignment,
#      which both changes the cells stored type and introduces binary float
#      noise (e.g. 418.26999999999998) and drops significant trailing zeros
#      (e.g. "1.00" -> 1). The source workbook keeps these as literal text, so
#      those specific cells are temporarily switched to Text number format,
#      written, then switched back via ClearFormats() so the on-disk style
#      index is unaffected (stays the default/unstyled index) while the
#      stored value remains textual.
#
#      The format toggle is batched through a single Union range (rather than
#      looping per cell) so only one extra style entry is ever created. The
#      engine only honours NumberFormat/ClearFormats on the FIRST area of a
#      multi-area Range when set on the union directly, so both calls iterate
#      `.Areas` explicitly to make sure every disjoint block is covered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Values that are safe to assign directly (stay text on their own) ---
$ws.Range("D2").Value = "65.424.25"
$ws.Range("E2").Value = "  +5.34%  "
$ws.Range("D3").Value = "3.508.79"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  +5.14%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +7.12%  "
$ws.Range("E10").Value = "  +17.41%  "
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("E12").Value = "  +23.28%  "
$ws.Range("E13").Value = "  +9.39%  "
$ws.Range("D14").Value = "4.062.48"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "3.520.38"
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D20").Value = "65.271.84"
$ws.Range("E20").Value = "  +5.03%  "
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("E27").Value = "  +2.84%  "
$ws.Range("E28").Value = "  +6.82%  "
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E31").Value = "  +5.80%  "
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("E36").Value = "  +4.34%  "
$ws.Range("D37").Value = "0.0₃0743"
$ws.Range("E37").Value = "  +38.35%  "
$ws.Range("E38").Value = "  +10.41%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  +5.17%  "
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  -2.79%  "
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("E49").Value = "  +5.45%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("E51").Value = "  +10.93%  "

# --- 2) Plain-decimal price values that must remain literal text ---
$numericTextValues = [ordered]@{
    "D5" = "418.27"
    "D6" = "132.56"
    "D7" = "0.656"
    "D8" = "0.999"
    "D9" = "0.778"
    "D10" = "0.164"
    "D11" = "43.47"
    "D12" = "0.0000268"
    "D13" = "10.05"
    "D15" = "0.141"
    "D16" = "20.58"
    "D18" = "12.81"
    "D21" = "455.43"
    "D22" = "90.39"
    "D24" = "13.38"
    "D26" = "9.96"
    "D27" = "34.22"
    "D28" = "12.63"
    "D30" = "7.45"
    "D31" = "0.118"
    "D33" = "39.82"
    "D34" = "1.00"
    "D36" = "0.0507"
    "D39" = "0.998"
    "D42" = "2.76"
    "D43" = "146.13"
    "D48" = "15.91"
    "D50" = "21.84"
    "D51" = "2.57"
}

$numericTextCells = @($numericTextValues.Keys)
$numericTextRange = $ws.Range($numericTextCells[0])
foreach ($addr in $numericTextCells[1..($numericTextCells.Length - 1)]) {
    $numericTextRange = $excel.Union($numericTextRange, $ws.Range($addr))
}

for ($i = 1; $i -le $numericTextRange.Areas.Count; $i++) {
    $numericTextRange.Areas.Item($i).NumberFormat = "@"
}
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Value = $numericTextValues[$addr]
}
for ($i = 1; $i -le $numericTextRange.Areas.Count; $i++) {
    $numericTextRange.Areas.Item($i).ClearFormats()
}

